$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '35.022.89'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.65%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.847.35'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.65%  '

# Row 4
$ws.Range('E4').Value = '  +0.18%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '232.09'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.38%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.618'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.83%  '

# Row 7
$ws.Range('E7').Value = '  +0.09%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '40.70'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.39%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.332'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.42%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0693'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.80%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0979'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.96%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.112.02'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.57%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.47'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +6.64%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.844.01'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.41%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.676'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.57%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.67'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.19%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '35.050.11'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.76%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '70.00'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.00%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0792'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.43%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '240.44'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.23%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.26'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.53%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.69'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.26%  '

# Row 23
$ws.Range('E23').Value = '  +0.18%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.27'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.98%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '172.44'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.90%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.84'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.46%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.56'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.44%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.124'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.26%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.55'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.21%  '

# Row 30
$ws.Range('E30').Value = '  +0.21%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0555'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.98%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.95'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.42%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.98'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.51%  '

# Row 34
$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.59'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +22.84%  '

# Row 35
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.96'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +12.71%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.757'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +13.02%  '

# Row 37
$ws.Range('E37').Value = '  +8.67%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.08'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +14.05%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '90.17'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.71%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.350.10'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.53%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0196'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.80%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '14.61'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.58%  '

# Row 43
$ws.Range('E43').Value = '  +4.64%  '

# Row 44
$ws.Range('E44').Value = '  -1.83%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.77'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.36%  '

# Row 46
$ws.Range('E46').Value = '  +4.32%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.31'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.41%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.029.47'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.71%  '

# Row 49
$ws.Range('B49').Value = 'THORChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.41'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +21.18%  '

# Row 50
$ws.Range('B50').Value = 'PaxDollar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.01'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.12%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0669'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.32%  '
